# Update the "Solar" column (E) values for years 2023 and 2024
# to incorporate updated data from upstream processes through 2024.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 25 -> Open year 2023, Solar (column E) 65 -> 66
$ws.Range("E25").Value = 66

# Row 26 -> Open year 2024, Solar (column E) 34 -> 48
$ws.Range("E26").Value = 48
